$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) "packages" sheet: bump the version/date in the results package
#    description (row 3, column C).
# ------------------------------------------------------------------
$packages = $wb.Worksheets.Item("packages")
$packages.Range("C3").Value2 = "Metadata on files produced by working group(v1.2.0, 2022-09-20)"

# ------------------------------------------------------------------
# 2) "entities" sheet: add two new working-group entities ("meta" and
#    "sandbox"), keeping the existing alphabetical ordering.
#    Columns: A=package B=name C=label D=description E=abstract F=extends
# ------------------------------------------------------------------
$entities = $wb.Worksheets.Item("entities")

# Fix pluralisation typo in the existing "denovo" row description.
$entities.Range("D6").Value2 = "Denovo working group files"

# Insert "meta" row right before the existing "proteomics" row (row 7).
$entities.Rows.Item(7).Insert()
$entities.Range("A7").Value2 = "rd3_cluster_results"
$entities.Range("B7").Value2 = "meta"
$entities.Range("C7").Value2 = "Meta Analysis"
$entities.Range("D7").Value2 = "Meta Analysis working group files"
$entities.Range("F7").Value2 = "rd3_cluster_results_template"

# Insert "sandbox" row right before the existing "snvindel" row, which is
# now row 10 after the previous insertion.
$entities.Rows.Item(10).Insert()
$entities.Range("A10").Value2 = "rd3_cluster_results"
$entities.Range("B10").Value2 = "sandbox"
$entities.Range("C10").Value2 = "Sandbox"
$entities.Range("D10").Value2 = "Sandbox only files"
$entities.Range("F10").Value2 = "rd3_cluster_results_template"

# ------------------------------------------------------------------
# 3) "attributes" sheet: rename the "xref" dataType to "categorical"
#    (used by the "ern" and "extension" attribute rows).
# ------------------------------------------------------------------
$attributes = $wb.Worksheets.Item("attributes")
$attributes.Range("D9").Value2 = "categorical"
$attributes.Range("D12").Value2 = "categorical"
